$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the Area / Atotal columns
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# Area formula for row 2 references a fixed 0 (no previous segment)
$ws.Range("G2").Formula = "=(D2-0)*B2/100"

# Area formula for row 3 (first real "previous row" reference, not yet shared)
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# Shared formula block G4:G15 following the same pattern as the D/E shared formulas
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Total area (sum of the individual segment areas)
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Update the current selection to match the edited workbook
$ws.Range("H2").Select()
